# Updated cryptos list on Tue Aug 13 13:36:13 UTC 2024 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# swaps the VeChain / WhiteBITCoin rows (49/50) to reflect their new order.
# Leading "'" forces plain decimal-looking prices to stay text (matching the
# original inlineStr cells) instead of being auto-coerced to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.847.77"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.634.34"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'520.03"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'146.32"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D9").Value = "2.639.47"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").Value = "'6.32"
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "3.094.65"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "58.836.49"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "'20.85"
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "2.638.71"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").Value = "'348.84"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "'4.47"
$ws.Range("E20").Value = "  -3.41%  "
$ws.Range("D21").Value = "'10.28"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").Value = "'6.17"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "'61.81"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("E25").Value = "  -2.80%  "
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").Value = "0.0₃0804"
$ws.Range("E28").Value = "  -3.58%  "
$ws.Range("D29").Value = "'7.07"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").Value = "'6.29"
$ws.Range("E31").Value = "  -5.06%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "'18.88"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").Value = "'149.30"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "'0.973"
$ws.Range("E35").Value = "  -6.66%  "
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").Value = "'36.52"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").Value = "'0.848"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("D40").Value = "'1.43"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").Value = "'280.72"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").Value = "'0.997"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").Value = "'0.0986"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("D45").Value = "'19.65"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("E46").Value = "  -4.69%  "
$ws.Range("D47").Value = "2.071.93"
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0230"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "'10.29"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").Value = "'4.68"
$ws.Range("E51").Value = "  -2.79%  "
